$p = $ppt.ActivePresentation

# --- Update the cached date field shown on the Notes Master footer -------
# (was "11/21/2019", becomes "8/24/2020")
$nm = $p.NotesMaster
$hf = $nm.HeadersFooters
$hf.DateAndTime.Text = "8/24/2020"

# --- Insert 3 new "section divider" slides --------------------------------
# Each is a duplicate of slide 1 ("Policy Instruments -3" / "Market power" /
# "Adverse selection" / "Enforcement"), inserted right before each of the
# three major topic sections.

# 1) After "0 Introduction" (slide 2), before "Price-taker" (slide 3)
$dup1 = $p.Slides.Item(1).Duplicate()
$dup1.Item(1).MoveTo(3)

# 2) After "Extensions" (now at position 12), before "Unknown abatement
#    costs"
$dup2 = $p.Slides.Item(1).Duplicate()
$dup2.Item(1).MoveTo(13)

# 3) After the second "Application" slide (now at position 20), before
#    "Enforcement"
$dup3 = $p.Slides.Item(1).Duplicate()
$dup3.Item(1).MoveTo(21)
